# New row of data (row 5), matching the FES05N layout of rows 2-4:
#   product code | product description | batch | qty | line no | remark

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A5, D5 and E5 hold digit-only text ("20001119", "1", "4"); force them to
# Text format first so Excel stores them as literal strings (shared-string
# table entries) instead of silently auto-converting to numbers.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("D5:E5").NumberFormat = "@"

$ws.Range("A5").Value = "20001119"
$ws.Range("B5").Value = "SOSRO TEH KOTAK 4+2S"
$ws.Range("C5").Value = "FES05N"
$ws.Range("D5").Value = "1"
$ws.Range("E5").Value = "4"
$ws.Range("F5").Value = "RT,(E-1B)"

# Match the bordered look of the rows above by copying row 4's formatting
# (border + General number format) onto the new row.
$ws.Range("A4:F4").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)
